$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45233 to 45243
# for every data row (rows 2 through 158).
$ws.Range("C2:C158").Value = 45243
